$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 95 (shifts existing rows 95:194 down to 96:195)
$ws.Rows("95:95").Insert()

# Populate the newly inserted row with the new Cilantro price record
$ws.Range("A95").Value = 8
$ws.Range("B95").Value = "Terminal La Palmera de La Serena"
$ws.Range("C95").Value = "Coquimbo"
$ws.Range("D95").Value = 44944
$ws.Range("E95").Value = 4
$ws.Range("F95").Value = 100112040
$ws.Range("G95").Value = "Cilantro"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 1600
$ws.Range("K95").Value = 3000
$ws.Range("L95").Value = 3500
$ws.Range("M95").Value = 3250
$ws.Range("N95").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O95").Value = "Provincia del Elquí"
$ws.Range("P95").Value = 2167
$ws.Range("Q95").Value = 1.5
$ws.Range("R95").Value = "Hortaliza"
